# Remove the post entry that was deleted from the spreadsheet
# (row 128: 「はじめてのアラビアもじ」). All rows below shift up by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(128).Delete()
